# Update of the DRomics TODO list - "short term" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("short term")

# Row 17 ("65. changer le nom adjpvalue en qvalue ...") - replace the follow-up
# paragraph: the team decided NOT to rename, but to mention it in the help.
$ws.Range("A17").Value = "65. changer le nom adjpvalue en qvalue dans les sorties pour que ce soit plus clair !!! (pas clair la tricherie dans Rager 2017)
Non car personne ne connait en ecotox, mais a mentionner dans l'aide"

# Row 18 ("69. mettre un warning ...") - fix typo (réplcats -> réplicats) and add a
# second sentence about forbidding anova when too many doses have a single replicate.
$ws.Range("A18").Value = "69. mettre un warning ou interdire l'utilisation de l'anova si pas au moins 2 ou 3 réplicats par groupe
interdire anova si trop de doses avec un seul réplicat (ex. superieur à 50%)"

# Person markers in column B
$ws.Range("B16").Value = "A"
$ws.Range("B18").Value = "ML"

# Row 18 now wraps onto two lines; match the natural two-line row height.
$ws.Rows.Item(18).RowHeight = 28.8

# Re-colour the status highlight: row 15 becomes yellow (like row 14), row 16
# becomes orange (like rows 17-18) to reflect the updated priority/status.
$ws.Range("A15").Interior.Color = 65535
$ws.Range("A16").Interior.Color = 49407

# Selection bookkeeping: the "done" sheet's last-selected cell moved too while
# reviewing, even though "short term" remains the active tab.
$wsDone = $wb.Worksheets.Item("done")
$wsDone.Range("A45").Select()

$ws.Select()
$ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 7
